$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("E6").Value = "  -4.32%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("E8").Value = "  -1.69%  "
$ws.Range("E9").Value = "  -1.59%  "
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("E12").Value = "  -3.86%  "
$ws.Range("E13").Value = "  -1.71%  "
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("E19").Value = "  -1.48%  "
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("E21").Value = "  -0.51%  "
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("E23").Value = "  -1.65%  "
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("E25").Value = "  -3.44%  "
$ws.Range("E26").Value = "  +1.94%  "
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("E28").Value = "  -2.75%  "
$ws.Range("E29").Value = "  +0.65%  "
$ws.Range("E31").Value = "  -1.51%  "
$ws.Range("E32").Value = "  +1.58%  "
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("E34").Value = "  -3.51%  "
$ws.Range("E35").Value = "  -7.42%  "
$ws.Range("E36").Value = "  -1.51%  "
$ws.Range("E37").Value = "  -4.30%  "
$ws.Range("E38").Value = "  -1.42%  "
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("E41").Value = "  -3.74%  "
$ws.Range("E42").Value = "  -7.43%  "
$ws.Range("E43").Value = "  -5.51%  "
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("E45").Value = "  +0.52%  "
$ws.Range("E46").Value = "  -2.21%  "
$ws.Range("E47").Value = "  +0.94%  "
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("E50").Value = "  -1.84%  "
$ws.Range("E51").Value = "  -2.87%  "
